$wb = $excel.ActiveWorkbook

# 2025 sheet: plain updated value
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 356760.80000000005

# 2030 sheet: formula referencing 2025
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Formula = "='2025'!B2*(1-0.2*0.2)"

# 2035 sheet: formula referencing 2025
$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("B2").Formula = "='2025'!B2*(1-0.2*0.4)"

# 2040 sheet: formula referencing 2025
$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("B2").Formula = "='2025'!B2*(1-0.2*0.6)"

# 2045 sheet: formula referencing 2025
$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("B2").Formula = "='2025'!B2*(1-0.2*0.8)"

# 2050 sheet: formula referencing 2025
$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("B2").Formula = "='2025'!B2*(1-0.2*1)"
